# Forward Look weekly refresh: 28 March 2024 -> 05 April 2024
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the "as at" date in the intro paragraph (A2)
$ws.Range("A2").Value = "This list contains a week-by-week view of  MoJ Official and National Statistics that have been pre-announced on the gov.uk release calendar as at 05 April 2024"

# 2. Flip specific publications from "provisional" to "confirmed" now their date has been verified
$confirmedRows = @(16, 17, 24, 27, 28, 29, 32, 35)
foreach ($r in $confirmedRows) {
    $ws.Range("D" + $r).Value = "confirmed"
}

# 3. Append new rows for weeks that rolled onto the bottom of the forward look (rows 78-80)
#    Start from the last existing data row (77) so the new rows inherit its style (s="5").
$ws.Range("A77:F77").Copy($ws.Range("A78:F80"))

# Row 78: week commencing only, no publication yet
$ws.Range("A78").Value = "03 Feb 2025"
$ws.Range("B78:D78").ClearContents()
$ws.Range("E78").Value = 6
$ws.Range("F78").ClearContents()

# Row 79: week commencing only, no publication yet
$ws.Range("A79").Value = "10 Feb 2025"
$ws.Range("B79:D79").ClearContents()
$ws.Range("E79").Value = 7
$ws.Range("F79").ClearContents()

# Row 80: new publication entry
$ws.Range("A80").Value = "17 Feb 2025"
$ws.Range("B80").Value = "HM Prison and Probation Service workforce quarterly: December 2024"
$ws.Range("C80").Value = "20 February 2025"
$ws.Range("D80").Value = "provisional"
$ws.Range("E80").Value = 8
$ws.Range("F80").Value = "standard"

# 4. Extend the conditional formatting ranges to cover the new rows (5:77 -> 5:80)
$addrFull = '$A$5:$F$77'
$addrFirstCol = '$A$5:$A$77'
$fcs = $ws.Range("A5:F77").FormatConditions
for ($i = 1; $i -le $fcs.Count(); $i++) {
    $fc = $fcs.Item($i)
    $addr = $fc.AppliesTo().Address()
    if ($addr -eq $addrFull) {
        $fc.ModifyAppliesToRange($ws.Range("A5:F80"))
    } elseif ($addr -eq $addrFirstCol) {
        $fc.ModifyAppliesToRange($ws.Range("A5:A80"))
    }
}
